# Added logic for age parameter in sorting
# Populate the "Age" column (column O) with per-row age values and apply
# the matching white-fill highlight that the authoring tool stamped on
# the newly-populated cells. Also nudge column P's stored width, which
# shifted slightly as a side effect of the same edit in the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ages = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 6
    6  = 6
    7  = 7
    8  = 5
    9  = 5
    10 = 7
    11 = 7
    12 = 7
    13 = 7
    14 = 6
    15 = 6
    16 = 7
    17 = 7
    18 = 6
    19 = 5
    20 = 6
    21 = 6
    22 = 5
    23 = 5
    24 = 6
    25 = 7
    26 = 7
    27 = 6
    28 = 5
    29 = 6
    30 = 6
    31 = 7
}

foreach ($row in $ages.Keys) {
    $ws.Cells.Item($row, 15).Value = $ages[$row]
}

# Column O (15) is "Age" - column 15, rows 2-31 carry the new values above.
$ws.Range("O2:O31").Interior.Color = 16777215

# Column P ("Requests", index 16) picked up a slightly wider stored width.
$ws.Columns.Item(16).ColumnWidth = 30.8
